$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 odds
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 2.75
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8
$ws.Range("O2").Value = 1.4
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65
$ws.Range("X2").Value = 9
$ws.Range("AN2").Value = 4

# Update row 3 odds
$ws.Range("J3").Value = 2.25
$ws.Range("L3").Value = 6.5
$ws.Range("O3").Value = 1.36
$ws.Range("P3").Value = 3
$ws.Range("Q3").Value = 2.2
$ws.Range("R3").Value = 1.65
$ws.Range("Z3").Value = 11
$ws.Range("AH3").Value = 13
$ws.Range("AJ3").Value = 21
$ws.Range("AZ3").Value = 151
$ws.Range("BA3").Value = 201

# Delete rows 4 and 5 (rows shift up automatically)
$ws.Rows("4:5").Delete()
